# Add data for 2022-09-12
#
# This updates the CTA violent-crime year-to-date workbook with a new
# batch of incidents recorded on 2022-09-12. The new incidents bump the
# cumulative YTD counts (and the corresponding "Total" rows/cells) on:
#   - the per-neighborhood worksheets (Robbery / Aggravated Assault /
#     Aggravated Battery rows, depending on the neighborhood),
#   - the citywide "Citywide Totals" rollup sheet, and
#   - the "By Neighborhood" rollup sheet.
# Jefferson Park gains a brand-new "Aggravated Battery" figure for 2019
# (previously empty / zero), so that cell is written for the first time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 50
$ws.Range("F3").Value = 95
$ws.Range("B6").Value = 280
$ws.Range("C6").Value = 348
$ws.Range("E6").Value = 319
$ws.Range("G6").Value = 357
$ws.Range("H6").Value = 324
$ws.Range("I6").Value = 383
$ws.Range("B7").Value = 381
$ws.Range("C7").Value = 470
$ws.Range("E7").Value = 478
$ws.Range("F7").Value = 557
$ws.Range("G7").Value = 518
$ws.Range("H7").Value = 514
$ws.Range("I7").Value = 635

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("E6").Value = 31
$ws.Range("E7").Value = 40

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("C6").Value = 25
$ws.Range("H6").Value = 20
$ws.Range("C7").Value = 30
$ws.Range("H7").Value = 32

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 8

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("B4").Value = 11
$ws.Range("B5").Value = 13

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("E6").Value = 17
$ws.Range("E7").Value = 24

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("C6").Value = 3
$ws.Range("H8").Value = 37
$ws.Range("E28").Value = 24
$ws.Range("E32").Value = 40
$ws.Range("C36").Value = 30
$ws.Range("H36").Value = 32
$ws.Range("F45").Value = 3
$ws.Range("G77").Value = 17
$ws.Range("I77").Value = 35
$ws.Range("B80").Value = 13
$ws.Range("G88").Value = 8
$ws.Range("B98").Value = 381
$ws.Range("C98").Value = 470
$ws.Range("E98").Value = 478
$ws.Range("F98").Value = 557
$ws.Range("G98").Value = 518
$ws.Range("H98").Value = 514
$ws.Range("I98").Value = 635

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("G6").Value = 12
$ws.Range("I6").Value = 22
$ws.Range("G7").Value = 17
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 3

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("H5").Value = 30
$ws.Range("H6").Value = 37

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("B2").Value = 2
$ws.Range("B5").Value = 3
